# edit.ps1 — applies the "completed a load of work not he project today" commit.
#
# Two independent edits to dissertation/project planning.docx:
#
#  1. The stray empty paragraph right after the "Software processes" heading
#     used to carry the document's "_GoBack" bookmark (an artifact Word drops
#     at the last edit position). That bookmark is removed, leaving a plain
#     empty paragraph.
#
#  2. A new closing discussion of the Akka/testKit unit-testing workaround is
#     appended after the "...so as to catch the most errors possible."
#     paragraph: a blank spacer paragraph followed by the new paragraph of
#     prose. The "_GoBack" bookmark now lands at the end of that new prose
#     paragraph (wherever the cursor was after the last bit of typing).

$d = $word.ActiveDocument

# --- 1. Drop the old _GoBack bookmark (paragraph right after "Software processes") ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- 2. Append the two new paragraphs after "...most errors possible." ---
$rng = $d.Content
$found = $rng.Find.Execute("most errors possible.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate the 'most errors possible.' paragraph"
}

$insertPoint = $rng.End
$target = $d.Range($insertPoint, $insertPoint)

$newParasXml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:jc w:val="both"/></w:pPr></w:p><w:p><w:pPr><w:jc w:val="both"/></w:pPr><w:r><w:t xml:space="preserve">Because I was using the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Akka</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> concurrency framework in order to make my application scalable. I could not use the standard Unit tests provided by Play framework to test all aspects of my application.  This is because the actors in play are protected by the special </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ActorRef</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> class, meaning that you can only send and receive messages to them</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>, .</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> The </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>akka</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>testKit</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> enables you </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>to  get</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> access to the underlying actors methods, meaning I could make full use of Spec2 to enable correct testing of my application. </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$target.InsertXML($newParasXml)

Write-Output "Applied project-planning edits."
